$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("loginData")
$ws.Activate()

# Update cell C4 value from "Invalid" to "Valid"
$ws.Range("C4").Value = "Valid"

# Update selection to C4 only
$ws.Range("C4").Select()

# Update window view properties
$wb.Windows.Item(1).Left = 2160
$wb.Windows.Item(1).Top = 2160
$wb.Windows.Item(1).Width = 14400
$wb.Windows.Item(1).Height = 7290
